# Updates the cryptos list (Coin/Link/Price/Volume(1h)) table on Sheet1.
# Refreshes Price (col D) and Volume(1h) % (col E) figures for the crypto rows,
# and swaps the InjectiveProtocol/Stellar rows (34/35) to reflect the new ranking
# order while refreshing their price/volume figures too.
#
# All of these source cells are plain text in the workbook (t="inlineStr"),
# including price strings that otherwise look numeric (e.g. "241.86") and
# percent strings padded with spaces (e.g. "  -0.50%  "). Assigning a plain
# numeric-looking string straight to Range.Value lets Excel's COM layer
# auto-convert it to a Number cell, so Set-TextValue forces text via
# NumberFormat "@" first, then restores the "Normal" style afterwards so no
# stray per-cell formatting (quote-prefix style) is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.013.41"
Set-TextValue $ws.Range("E2") "  -0.50%  "
Set-TextValue $ws.Range("D3") "2.215.55"
Set-TextValue $ws.Range("E3") "  -1.36%  "
Set-TextValue $ws.Range("E4") "  +0.22%  "
Set-TextValue $ws.Range("D5") "241.86"
Set-TextValue $ws.Range("E5") "  -1.87%  "
Set-TextValue $ws.Range("D6") "0.628"
Set-TextValue $ws.Range("E6") "  -0.45%  "
Set-TextValue $ws.Range("D7") "73.09"
Set-TextValue $ws.Range("E7") "  -1.90%  "
Set-TextValue $ws.Range("E8") "  +0.20%  "
Set-TextValue $ws.Range("D9") "0.606"
Set-TextValue $ws.Range("E9") "  -2.16%  "
Set-TextValue $ws.Range("D10") "42.27"
Set-TextValue $ws.Range("E10") "  -0.61%  "
Set-TextValue $ws.Range("D11") "0.0954"
Set-TextValue $ws.Range("E11") "  +0.57%  "
Set-TextValue $ws.Range("E12") "  -1.58%  "
Set-TextValue $ws.Range("D13") "0.104"
Set-TextValue $ws.Range("E13") "  +0.48%  "
Set-TextValue $ws.Range("D14") "2.549.15"
Set-TextValue $ws.Range("E14") "  -1.20%  "
Set-TextValue $ws.Range("D15") "14.27"
Set-TextValue $ws.Range("E15") "  -1.47%  "
Set-TextValue $ws.Range("D16") "0.836"
Set-TextValue $ws.Range("E16") "  -1.75%  "
Set-TextValue $ws.Range("D17") "2.207.50"
Set-TextValue $ws.Range("E17") "  -2.19%  "
Set-TextValue $ws.Range("D18") "41.910.29"
Set-TextValue $ws.Range("E18") "  -0.48%  "
Set-TextValue $ws.Range("E19") "  +7.56%  "
Set-TextValue $ws.Range("D20") "6.20"
Set-TextValue $ws.Range("E20") "  +1.16%  "
Set-TextValue $ws.Range("D21") "72.98"
Set-TextValue $ws.Range("E21") "  +0.84%  "
Set-TextValue $ws.Range("D22") "10.54"
Set-TextValue $ws.Range("E22") "  +18.08%  "
Set-TextValue $ws.Range("D23") "230.30"
Set-TextValue $ws.Range("E23") "  -0.28%  "
Set-TextValue $ws.Range("E24") "  -6.68%  "
Set-TextValue $ws.Range("E25") "  +4.07%  "
Set-TextValue $ws.Range("E26") "  +0.10%  "
Set-TextValue $ws.Range("D27") "3.67"
Set-TextValue $ws.Range("E27") "  +1.37%  "
Set-TextValue $ws.Range("E28") "  -1.28%  "
Set-TextValue $ws.Range("E29") "  -2.84%  "
Set-TextValue $ws.Range("D30") "168.43"
Set-TextValue $ws.Range("E30") "  -0.54%  "
Set-TextValue $ws.Range("E31") "  -0.81%  "
Set-TextValue $ws.Range("D32") "5.62"
Set-TextValue $ws.Range("E32") "  +7.13%  "
Set-TextValue $ws.Range("D33") "0.0796"
Set-TextValue $ws.Range("E33") "  -2.75%  "
Set-TextValue $ws.Range("B34") "InjectiveProtocol"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D34") "29.70"
Set-TextValue $ws.Range("E34") "  -4.25%  "
Set-TextValue $ws.Range("B35") "Stellar"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D35") "0.125"
Set-TextValue $ws.Range("E35") "  -0.21%  "
Set-TextValue $ws.Range("E36") "  -9.73%  "
Set-TextValue $ws.Range("E37") "  -3.77%  "
Set-TextValue $ws.Range("D38") "0.0301"
Set-TextValue $ws.Range("E38") "  -4.17%  "
Set-TextValue $ws.Range("D39") "13.68"
Set-TextValue $ws.Range("E39") "  -1.09%  "
Set-TextValue $ws.Range("D40") "66.00"
Set-TextValue $ws.Range("E40") "  +5.09%  "
Set-TextValue $ws.Range("E41") "  -2.15%  "
Set-TextValue $ws.Range("E42") "  -2.29%  "
Set-TextValue $ws.Range("E43") "  -2.97%  "
Set-TextValue $ws.Range("D44") "8.80"
Set-TextValue $ws.Range("E44") "  +1.28%  "
Set-TextValue $ws.Range("D45") "105.08"
Set-TextValue $ws.Range("E45") "  -1.70%  "
Set-TextValue $ws.Range("E46") "  -2.22%  "
Set-TextValue $ws.Range("D47") "2.43"
Set-TextValue $ws.Range("E47") "  +5.16%  "
Set-TextValue $ws.Range("E48") "  -0.19%  "
Set-TextValue $ws.Range("E49") "  -0.43%  "
Set-TextValue $ws.Range("E50") "  -0.18%  "
Set-TextValue $ws.Range("D51") "2.424.32"
Set-TextValue $ws.Range("E51") "  -1.43%  "
